# "Adding test cases for watch lists"
#
# - Resize the workbook's application window (cosmetic bookViews entry).
# - Add a new "SKIP" result value (becomes a new shared string) and apply it
#   to the Results column (E) for the first three test cases (rows 2-4),
#   replacing "PASS".
# - Move the sheet selection from C7 to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Resize the Excel window (reflected in xl/workbook.xml bookViews/workbookView).
$win = $excel.ActiveWindow
$win.Width = 15150
$win.Height = 10125

# Mark the first three watch-list test cases (E2:E4) as "SKIP" instead of "PASS".
$ws.Range("E2").Value = "SKIP"
$ws.Range("E3").Value = "SKIP"
$ws.Range("E4").Value = "SKIP"

# Update the active selection on the sheet from C7 to C6.
$ws.Range("C6").Select() | Out-Null
